$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.276772
$ws.Range("H2").Value = 12.830316
$ws.Range("I2").Value = 0.06135676581847978
$ws.Range("J2").Value = 0.06135676581847978
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 117.044563
$ws.Range("N2").Value = 351.133689
$ws.Range("O2").Value = 0.3245365645427815
$ws.Range("P2").Value = 0.3245365645427815
$ws.Range("Q2").Value = 500.572909790636
$ws.Range("R2").Value = 4505.156188115724
$ws.Range("S2").Value = 0.01991251399018539
$ws.Range("T2").Value = 0.01991251399018539

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.276772
$ws.Range("H3").Value = 12.830316
$ws.Range("I3").Value = 0.06135676581847978
$ws.Range("J3").Value = 0.06135676581847978
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.281657135515876
$ws.Range("P3").Value = 0.281657135515876
$ws.Range("Q3").Value = 434.4346594261547
$ws.Range("R3").Value = 3909.911934835392
$ws.Range("S3").Value = 0.01728157090495143
$ws.Range("T3").Value = 0.01728157090495143

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.276772
$ws.Range("H4").Value = 12.830316
$ws.Range("I4").Value = 0.06135676581847978
$ws.Range("J4").Value = 0.06135676581847978
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.3938062999413425
$ws.Range("P4").Value = 0.3938062999413425
$ws.Range("Q4").Value = 607.4161958706987
$ws.Range("R4").Value = 5466.745762836288
$ws.Range("S4").Value = 0.02416268092334296
$ws.Range("T4").Value = 0.02416268092334296

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 36.44531133333333
$ws.Range("H5").Value = 109.335934
$ws.Range("I5").Value = 0.5228631389891535
$ws.Range("J5").Value = 0.5228631389891535
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 117.044563
$ws.Range("N5").Value = 351.133689
$ws.Range("O5").Value = 0.3245365645427815
$ws.Range("P5").Value = 0.3245365645427815
$ws.Range("Q5").Value = 4265.725538408947
$ws.Range("R5").Value = 38391.52984568052
$ws.Range("S5").Value = 0.1696882068535947
$ws.Range("T5").Value = 0.1696882068535947

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 36.44531133333333
$ws.Range("H6").Value = 109.335934
$ws.Range("I6").Value = 0.5228631389891535
$ws.Range("J6").Value = 0.5228631389891535
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.281657135515876
$ws.Range("P6").Value = 0.281657135515876
$ws.Range("Q6").Value = 3702.116085864956
$ws.Range("R6").Value = 33319.0447727846
$ws.Range("S6").Value = 0.1472681339945243
$ws.Range("T6").Value = 0.1472681339945243

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 36.44531133333333
$ws.Range("H7").Value = 109.335934
$ws.Range("I7").Value = 0.5228631389891535
$ws.Range("J7").Value = 0.5228631389891535
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.3938062999413425
$ws.Range("P7").Value = 0.3938062999413425
$ws.Range("Q7").Value = 5176.210554927079
$ws.Range("R7").Value = 46585.8949943437
$ws.Range("S7").Value = 0.2059067981410344
$ws.Range("T7").Value = 0.2059067981410344

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 28.98126466666666
$ws.Range("H8").Value = 86.943794
$ws.Range("I8").Value = 0.4157800951923667
$ws.Range("J8").Value = 0.4157800951923668
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 117.044563
$ws.Range("N8").Value = 351.133689
$ws.Range("O8").Value = 0.3245365645427815
$ws.Range("P8").Value = 0.3245365645427815
$ws.Range("Q8").Value = 3392.09945809734
$ws.Range("R8").Value = 30528.89512287607
$ws.Range("S8").Value = 0.1349358436990014
$ws.Range("T8").Value = 0.1349358436990014

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 28.98126466666666
$ws.Range("H9").Value = 86.943794
$ws.Range("I9").Value = 0.4157800951923667
$ws.Range("J9").Value = 0.4157800951923668
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 101.5800373333333
$ws.Range("N9").Value = 304.740112
$ws.Range("O9").Value = 0.281657135515876
$ws.Range("P9").Value = 0.281657135515876
$ws.Range("Q9").Value = 2943.917946807214
$ws.Range("R9").Value = 26495.26152126493
$ws.Range("S9").Value = 0.1171074306164003
$ws.Range("T9").Value = 0.1171074306164003

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 28.98126466666666
$ws.Range("H10").Value = 86.943794
$ws.Range("I10").Value = 0.4157800951923667
$ws.Range("J10").Value = 0.4157800951923668
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 142.0267893333333
$ws.Range("N10").Value = 426.080368
$ws.Range("O10").Value = 0.3938062999413425
$ws.Range("P10").Value = 0.3938062999413425
$ws.Range("Q10").Value = 4116.115971426243
$ws.Range("R10").Value = 37045.04374283619
$ws.Range("S10").Value = 0.1637368208769651
$ws.Range("T10").Value = 0.1637368208769651

